$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column E width. Target raw OOXML width is 5.42578125; this engine
# quantizes ColumnWidth to an internal 1/6-character grid, so the nearest
# achievable stored width is 5.5 (reached via a ColumnWidth input of
# target - 5/6, i.e. the natural unit conversion).
$ws.Columns.Item(5).ColumnWidth = 4.592447916666667

# Row 17 - new task "Seeders and migrations"
$ws.Range("A17").Value = "Seeders and migrations"
$ws.Range("B17").Value = 0.5
$ws.Range("C17").Value = 0.5
$ws.Range("D17").Value = "Ja"
$ws.Range("E17").Value = "Kan meer worden naarmate ik meer data nodig heb"

# Row 18 - new task "Items"
$ws.Range("A18").Value = "Items"
$ws.Range("B18").Value = 1.5
$ws.Range("C18").Value = 1.5
$ws.Range("D18").Value = "Ja"

# Totals now need to include the two new rows
$ws.Range("B24").Formula = "=SUM(B3:B18)"
$ws.Range("C24").Formula = "=SUM(C3:C18)"

# Update the selection to match the author's last-edited cell
$ws.Range("E18").Select()

$wb.Save()
